$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 38/39 swap: RenderToken <-> VeChain (rank order changed) ---
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01959"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.79%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.436"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.96%  "

# --- Price (column D) updates ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.907.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.810.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4627"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3757"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07480"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8799"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.45"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.778.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.360"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.548"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07049"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008759"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.910.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.316"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.999.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.924"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.152"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.303"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08904"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7703"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.168"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.483"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.902"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.118"
$ws.Range("D37").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05241"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5331"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.225"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.910"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1661"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.591"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5050"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.672"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.000"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06330"
$ws.Range("D51").Style = "Normal"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  -1.67%  "
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  +3.57%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("E10").Value = "  -0.69%  "
$ws.Range("E11").Value = "  -2.63%  "
$ws.Range("E12").Value = "  -2.57%  "
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("E14").Value = "  -3.03%  "
$ws.Range("E16").Value = "  -2.67%  "
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("E20").Value = "  -3.10%  "
$ws.Range("E21").Value = "  -1.69%  "
$ws.Range("E22").Value = "  +1.12%  "
$ws.Range("E23").Value = "  -0.81%  "
$ws.Range("E24").Value = "  -2.78%  "
$ws.Range("E25").Value = "  -1.87%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("E28").Value = "  -9.58%  "
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("E30").Value = "  -1.92%  "
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("E32").Value = "  -1.77%  "
$ws.Range("E33").Value = "  -2.35%  "
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("E35").Value = "  -0.98%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("E40").Value = "  -1.79%  "
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("E43").Value = "  +1.92%  "
$ws.Range("E44").Value = "  -3.67%  "
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("E47").Value = "  -2.83%  "
$ws.Range("E48").Value = "  -1.11%  "
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("E51").Value = "  -0.77%  "
